# Apply edit: "cam binding to camkii IS bimolecular, previous commit was wrong"
# Adds parameter values for CaMKII T286 and T306 phosphorylation rows (C10, C11),
# highlights the newly added blank reference cells (F8:G9), and moves selection to C10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New numeric values for the two phosphorylation rate parameters.
$ws.Range("C10").Value = 10
$ws.Range("C11").Value = 0.02

# Highlight the newly-present (previously absent) reference cells with the
# same highlight fill used to flag missing info.
$highlight = $ws.Range("F8:G9")
$highlight.Interior.ThemeColor = 9
$highlight.Interior.TintAndShade = 0.59999389629810485

# Move the active selection to the cell that was actually edited.
$ws.Range("C10").Select()
